$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Append a new "GuilID" field row (row 11) to the Property struct sheet,
# describing the guild's id ("工会ID") on the object-typed "Friend" relation.
$row = 11
$ws.Cells.Item($row, 1).Value = "GuilID"
$ws.Cells.Item($row, 2).Value = "object"
$ws.Cells.Item($row, 3).Value = $true
$ws.Cells.Item($row, 4).Value = $true
$ws.Cells.Item($row, 5).Value = $true
$ws.Cells.Item($row, 6).Value = $true
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = "Friend"
$ws.Cells.Item($row, 10).Value = "工会ID"

# Match the text-formatted number format ("@", style index 1) used by the
# Id/Type/RelationValue/Desc columns on the other data rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 10).NumberFormat = "@"

$ws.Range("E19").Select()
